$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 55: C 12 -> 13, E recalculated (C/D)
$ws.Range("C55").Value = 13
$ws.Range("E55").Value = 0.09285714285714286

# Row 113: C 18 -> 20, E recalculated
$ws.Range("C113").Value = 20
$ws.Range("E113").Value = 0.2298850574712644

# Row 181: C 22 -> 23, E recalculated
$ws.Range("C181").Value = 23
$ws.Range("E181").Value = 0.4893617021276596

# Row 185: C 37 -> 38, E recalculated
$ws.Range("C185").Value = 38
$ws.Range("E185").Value = 0.4871794871794872

# Row 188: C 14 -> 15, E recalculated
$ws.Range("C188").Value = 15
$ws.Range("E188").Value = 0.3488372093023256

# Row 190: C 18 -> 20, E recalculated
$ws.Range("C190").Value = 20
$ws.Range("E190").Value = 0.5263157894736842

# Row 191: C 16 -> 20, D 16 -> 20, E stays 1
$ws.Range("C191").Value = 20
$ws.Range("D191").Value = 20
$ws.Range("E191").Value = 1
